$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Runsheet update: add "exclude " / "reason_excluded" columns (I/J) so the
# demographics sheet records which participants were excluded and why.
#
# Rows 2-46 (subjects 1-45) get an inclusion flag in column I and a reason
# in column J. Three subjects are marked for exclusion:
#   row 13 (subject 12) - could not calibrate
#   row 24 (subject 23) - didn't calibrate
#   row 33 (subject 32) - started crying / fussed out
# Everybody else in that range is "include" / "none".
#
# The write order below matters: new shared-string entries are appended in
# first-use order, and this sequence reproduces the same order the runsheet
# ended up with (include, exclude, "exclude " header, reason_excluded
# header, none, calibration failure, fussed out).

$ws.Range("I2:I46").Value = "include"
$ws.Range("I13").Value = "exclude"
$ws.Range("I24").Value = "exclude"
$ws.Range("I33").Value = "exclude"
$ws.Range("I1").Value = "exclude "

$ws.Range("J1").Value = "reason_excluded"

$ws.Range("J2:J46").Value = "none"
$ws.Range("J13").Value = "calibration failure"
$ws.Range("J24").Value = "calibration failure"
$ws.Range("J33").Value = "fussed out"

[void]$ws.Range("J34:J46").Select()
